# atualizacao 2a RQ 2024-2028
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new agent/name rows after the existing last row (114).
# New distinct strings are introduced in this order so the shared-string
# table grows in the same sequence as the target workbook:
#   LIGHT SESA, ENF, EBO, Ceraçá, CERAL-DIS, PACTO
$ws.Range("A115").Value = "LIGHT SESA"
$ws.Range("B115").Value = "LIGHT"

$ws.Range("A116").Value = "EQUATORIAL GO"
$ws.Range("B116").Value = "EQUATORIAL GO"

$ws.Range("A117").Value = "ENF"
$ws.Range("B117").Value = "EMR"

$ws.Range("A118").Value = "EBO"
$ws.Range("B118").Value = "EPB"

$ws.Range("A119").Value = "Ceraçá"
$ws.Range("B119").Value = "OUTRA"

$ws.Range("A120").Value = "CERAL-DIS"
$ws.Range("B120").Value = "OUTRA"

$ws.Range("A121").Value = "EAC"
$ws.Range("B121").Value = "EAC"

# Update existing mappings: FORCEL -> PACTO for two agent rows
$ws.Range("B97").Value = "PACTO"
$ws.Range("B113").Value = "PACTO"

# Update selection to match the saved view state in the diff
$ws.Range("A97").Select()
